{"js": "// 1) Fix the filename reference in the first line of the document:\n//    \"-Date-02-27-2024-Day-29.docx\" -> \"-Date-02-27-2024-Day-28.docx\"\nconst dayResults = context.document.body.search(\"Day-29\", { matchCase: true, matchWholeWord: false });\ndayResults.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < dayResults.items.length; i++) {\n  dayResults.items[i].insertText(\"Day-28\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2) Collapse the two-run spelling of \"True\" (\"T\" + \"rue\") into a single\n//    run containing \"True\" (happens twice in the document, for Case 1 and\n//    Case 2 sample output). A case-sensitive, whole-word search for \"True\"\n//    only matches these two spots (the rest of the document only contains\n//    the lowercase word \"true\").\nconst trueResults = context.document.body.search(\"True\", { matchCase: true, matchWholeWord: true });\ntrueResults.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < trueResults.items.length; i++) {\n  trueResults.items[i].insertText(\"True\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Fix the filename reference in the first line of the document:\n#    \"-Date-02-27-2024-Day-29.docx\" -> \"-Date-02-27-2024-Day-28.docx\"\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"Day-29\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"Day-28\"\n$find.Execute([ref]\"Day-29\", $true, $false, $false, $false, $false, $true, 1, $false, \"Day-28\", 2)\n\n# 2) Collapse the two-run spelling of \"True\" (\"T\" + \"rue\") into a single\n#    run containing \"True\" (happens twice in the document, for Case 1 and\n#    Case 2 sample output). A case-sensitive, whole-word search for \"True\"\n#    only matches these two spots (the rest of the document only contains\n#    the lowercase word \"true\").\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Text = \"True\"\n$find2.MatchCase = $true\n$find2.MatchWholeWord = $true\n$find2.Replacement.ClearFormatting()\n$find2.Replacement.Text = \"True\"\n$find2.Execute([ref]\"True\", $true, $true, $false, $false, $false, $true, 1, $false, \"True\", 2)\n"}
